$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("psicon")

# Add new worksheet "Sheet1" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Sheet1"

# Copy column A values from psicon into Sheet1 (A1:A151)
for ($r = 1; $r -le 151; $r++) {
    $val = $ws1.Cells.Item($r, 1).Value2
    $ws2.Cells.Item($r, 1).Value = $val
}

# Reproduce the formatting of A1 (bold, centered/top, thin border) on the new sheet
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws2.Range("A1").Value = 0

# Header + failed rows, written in an order that reproduces the target
# shared-string table order: Status(150), FAILED(151), Reason(152), reason-text(153)
$reason = "matlab.engine.EngineError: Unable to launch MVM server: License Error: Error checking out licens"
$failedRows = @(2, 37, 52, 139, 140)

$ws2.Cells.Item(1, 2).Value = "Status"
foreach ($r in $failedRows) {
    $ws2.Cells.Item($r, 2).Value = "FAILED"
}

$ws2.Cells.Item(1, 3).Value = "Reason"
foreach ($r in $failedRows) {
    $ws2.Cells.Item($r, 3).Value = $reason
}

# Column C width (~23.33 characters, matching the target column width)
$ws2.Columns.Item(3).ColumnWidth = 22.5

# Selection/view settings
$ws1.Activate()
$ws1.Range("A1:A1048576").Select()

$ws2.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws2.Range("C141").Select()
